# big update, app is now being developed on mobile device
# -mobile controls implemented
# -database loading on mobile device has been fixed
#
# Spreadsheet side of the commit: mark two Furniture rows as complete
# (Author/Complete? columns) and add two new Furniture entries
# (Security Table, Chair) right after the existing "Rack" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark existing rows as finished -----------------------------------
# Row 7  : Hospital bed      -> Complete? = Y
# Row 8  : weapons table     -> Author = Sasha, Complete? = Y
# (write these first so the new "Y" shared string is interned before the
#  strings introduced by the new rows below, matching save order)
$ws.Range("D7").Value = "Y"
$ws.Range("C8").Value = "Sasha"
$ws.Range("D8").Value = "Y"

# --- Insert two new Furniture rows between "Rack" (row 9) and
#     "Weapon rack" (row 10), pushing everything below down by two ------
$ws.Rows("10:11").Insert()

$ws.Range("A10").Value = "Security Table"
$ws.Range("B10").Value = "a table with monitors, for viewing security camera footage "

$ws.Range("A11").Value = "Chair"
$ws.Range("B11").Value = "a chair for tables "

# --- Restore the active selection left behind by the edit ---------------
$ws.Range("D15").Select()
